$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing rows 2-7 (re-run timestamps + recomputed KPI metrics) ----
$ws.Range("B2").Value = 44715.57863613426
$ws.Range("T2").Value = 16.3481348028805
$ws.Range("U2").Value = 0.8294041735018478
$ws.Range("V2").Value = 0.6396518258887574
$ws.Range("W2").Value = 6700.0
$ws.Range("X2").Value = 214.71
$ws.Range("Y2").Value = 53.91804380215612
$ws.Range("B3").Value = 44715.578636342594
$ws.Range("T3").Value = 15.653157735820223
$ws.Range("U3").Value = 1.150769170647169
$ws.Range("V3").Value = 0.2739653785483065
$ws.Range("W3").Value = 7900.0
$ws.Range("X3").Value = 221.71
$ws.Range("Y3").Value = 81.8781765121124
$ws.Range("B4").Value = 44715.57863640046
$ws.Range("T4").Value = 23.39216620388452
$ws.Range("U4").Value = 1.607201052098019
$ws.Range("W4").Value = 5600.0
$ws.Range("X4").Value = 131.3
$ws.Range("B5").Value = 44715.57863646991
$ws.Range("T5").Value = 30.032700517876982
$ws.Range("U5").Value = 2.956545957017598
$ws.Range("W5").Value = 6600.0
$ws.Range("X5").Value = 94.18
$ws.Range("B6").Value = 44715.57863649305
$ws.Range("AH6").Value = 147.13026571798542
$ws.Range("AI6").Value = 40.59442340173604
$ws.Range("AJ6").Value = 57.28768982389991
$ws.Range("B7").Value = 44715.578636516206
$ws.Range("AL7").Value = 33.0949220463136
$ws.Range("AM7").Value = 44032.0

# ---- Append new run (run id = 2): rows 8-13 ----
# row 8
$ws.Cells.Item(8, 1).Value = 2.0
$ws.Cells.Item(8, 2).Value = 44715.58029042824
$ws.Cells.Item(8, 3).Value = 0.0
$ws.Cells.Item(8, 4).Value = "pop_buurten[0]"
$ws.Cells.Item(8, 5).Value = "rijtjeshuizen (laag)"
$ws.Cells.Item(8, 6).Value = 40.0
$ws.Cells.Item(8, 7).Value = 40.0
$ws.Cells.Item(8, 8).Value = 7.6923076923076925
$ws.Cells.Item(8, 9).Value = 15.0
$ws.Cells.Item(8, 10).Value = 0.0
$ws.Cells.Item(8, 11).Value = 27.5
$ws.Cells.Item(8, 12).Value = 55.0
$ws.Cells.Item(8, 13).Value = 2.5
$ws.Cells.Item(8, 14).Value = 0.0
$ws.Cells.Item(8, 15).Value = 0.1
$ws.Cells.Item(8, 16).Value = 95.0
$ws.Cells.Item(8, 17).Value = 0.0
$ws.Cells.Item(8, 18).Value = 5.0
$ws.Cells.Item(8, 19).Value = 0.0
$ws.Cells.Item(8, 20).Value = 16.3481348028805
$ws.Cells.Item(8, 21).Value = 0.8294041735018478
$ws.Cells.Item(8, 22).Value = 0.6396518258887574
$ws.Cells.Item(8, 23).Value = 6700.0
$ws.Cells.Item(8, 24).Value = 214.71
$ws.Cells.Item(8, 25).Value = 53.91804380215612
$ws.Cells.Item(8, 26).Value = 0.09904138006205263
$ws.Cells.Item(8, 27).Value = -0.0005016326660606007
$ws.Cells.Item(8, 28).Value = 0.0

# row 9
$ws.Cells.Item(9, 1).Value = 2.0
$ws.Cells.Item(9, 2).Value = 44715.58029063657
$ws.Cells.Item(9, 3).Value = 1.0
$ws.Cells.Item(9, 4).Value = "pop_buurten[1]"
$ws.Cells.Item(9, 5).Value = "VVD-wijk (laag)"
$ws.Cells.Item(9, 6).Value = 40.0
$ws.Cells.Item(9, 7).Value = 22.5
$ws.Cells.Item(9, 8).Value = 30.0
$ws.Cells.Item(9, 9).Value = 62.5
$ws.Cells.Item(9, 10).Value = 22.5
$ws.Cells.Item(9, 11).Value = 0.0
$ws.Cells.Item(9, 12).Value = 0.0
$ws.Cells.Item(9, 13).Value = 15.0
$ws.Cells.Item(9, 14).Value = 0.0
$ws.Cells.Item(9, 15).Value = 0.1
$ws.Cells.Item(9, 16).Value = 92.5
$ws.Cells.Item(9, 17).Value = 0.0
$ws.Cells.Item(9, 18).Value = 7.5
$ws.Cells.Item(9, 19).Value = 0.0
$ws.Cells.Item(9, 20).Value = 15.653157735820223
$ws.Cells.Item(9, 21).Value = 1.150769170647169
$ws.Cells.Item(9, 22).Value = 0.2739653785483065
$ws.Cells.Item(9, 23).Value = 7900.0
$ws.Cells.Item(9, 24).Value = 221.71
$ws.Cells.Item(9, 25).Value = 81.8781765121124
$ws.Cells.Item(9, 26).Value = 0.39611086602817985
$ws.Cells.Item(9, 27).Value = -0.0020463473673181807
$ws.Cells.Item(9, 28).Value = 0.0

# row 10
$ws.Cells.Item(10, 1).Value = 2.0
$ws.Cells.Item(10, 2).Value = 44715.580290694445
$ws.Cells.Item(10, 3).Value = 2.0
$ws.Cells.Item(10, 4).Value = "pop_buurten[2]"
$ws.Cells.Item(10, 5).Value = "rijtjeshuizen (hoog)"
$ws.Cells.Item(10, 6).Value = 40.0
$ws.Cells.Item(10, 7).Value = 87.5
$ws.Cells.Item(10, 8).Value = 42.5531914893617
$ws.Cells.Item(10, 9).Value = 35.0
$ws.Cells.Item(10, 10).Value = 0.0
$ws.Cells.Item(10, 11).Value = 12.5
$ws.Cells.Item(10, 12).Value = 52.5
$ws.Cells.Item(10, 13).Value = 0.0
$ws.Cells.Item(10, 14).Value = 0.0
$ws.Cells.Item(10, 15).Value = 0.225
$ws.Cells.Item(10, 16).Value = 0.0
$ws.Cells.Item(10, 17).Value = 0.0
$ws.Cells.Item(10, 18).Value = 12.5
$ws.Cells.Item(10, 19).Value = 87.5
$ws.Cells.Item(10, 20).Value = 23.39216620388452
$ws.Cells.Item(10, 21).Value = 1.607201052098019
$ws.Cells.Item(10, 22).Value = -0.0
$ws.Cells.Item(10, 23).Value = 5600.0
$ws.Cells.Item(10, 24).Value = 131.3
$ws.Cells.Item(10, 25).Value = 100.0
$ws.Cells.Item(10, 26).Value = -0.17531039169365872
$ws.Cells.Item(10, 27).Value = -0.009039148214032521
$ws.Cells.Item(10, 28).Value = 0.0

# row 11
$ws.Cells.Item(11, 1).Value = 2.0
$ws.Cells.Item(11, 2).Value = 44715.58029074074
$ws.Cells.Item(11, 3).Value = 3.0
$ws.Cells.Item(11, 4).Value = "pop_buurten[3]"
$ws.Cells.Item(11, 5).Value = "VVD-wijk (hoog)"
$ws.Cells.Item(11, 6).Value = 40.0
$ws.Cells.Item(11, 7).Value = 77.5
$ws.Cells.Item(11, 8).Value = 100.0
$ws.Cells.Item(11, 9).Value = 72.5
$ws.Cells.Item(11, 10).Value = 22.5
$ws.Cells.Item(11, 11).Value = 0.0
$ws.Cells.Item(11, 12).Value = 0.0
$ws.Cells.Item(11, 13).Value = 5.0
$ws.Cells.Item(11, 14).Value = 0.0
$ws.Cells.Item(11, 15).Value = 0.15
$ws.Cells.Item(11, 16).Value = 0.0
$ws.Cells.Item(11, 17).Value = 0.0
$ws.Cells.Item(11, 18).Value = 15.0
$ws.Cells.Item(11, 19).Value = 85.0
$ws.Cells.Item(11, 20).Value = 30.032700517876982
$ws.Cells.Item(11, 21).Value = 2.956545957017598
$ws.Cells.Item(11, 22).Value = -0.0
$ws.Cells.Item(11, 23).Value = 6600.0
$ws.Cells.Item(11, 24).Value = 94.18
$ws.Cells.Item(11, 25).Value = 100.0
$ws.Cells.Item(11, 26).Value = -0.7120066804492732
$ws.Cells.Item(11, 27).Value = -0.0027013259258631857
$ws.Cells.Item(11, 28).Value = 0.0

# row 12
$ws.Cells.Item(12, 1).Value = 2.0
$ws.Cells.Item(12, 2).Value = 44715.580290763886
$ws.Cells.Item(12, 29).Value = 0.0
$ws.Cells.Item(12, 30).Value = "holonAgent[0]"
$ws.Cells.Item(12, 31).Value = "WindHolon"
$ws.Cells.Item(12, 32).Value = 70.0
$ws.Cells.Item(12, 33).Value = 69.0
$ws.Cells.Item(12, 34).Value = 147.13026571798542
$ws.Cells.Item(12, 35).Value = 40.59442340173604
$ws.Cells.Item(12, 36).Value = 57.28768982389991

# row 13
$ws.Cells.Item(13, 1).Value = 2.0
$ws.Cells.Item(13, 2).Value = 44715.58029078704
$ws.Cells.Item(13, 37).Value = "root"
$ws.Cells.Item(13, 38).Value = 33.0949220463136
$ws.Cells.Item(13, 39).Value = 44032.0

# Apply the same date number-format (yyyy-mm-dd) used by column B elsewhere in the sheet
for ($r = 8; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).NumberFormat = "yyyy-mm-dd"
}
